# Update "Pais" sheet: refresh timestamp, swap four pairs of country
# labels (source data re-sorted by rank) and update the covid statistics
# columns (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) for the rows whose numbers moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Title / last-updated timestamp (row 1) ---------------------------
$ws.Range("A1").Value = "Datos actualizados a 16 de Septiembre de 2020 a las 14:35"

# --- Country label swaps (column A) ------------------------------------
$ws.Range("A43").Value  = "Paises Bajos"
$ws.Range("A44").Value  = "China"

$ws.Range("A114").Value = "Uganda"
$ws.Range("A115").Value = "Suazilandia"

$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("A205").Value = "Santa Lucia"

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("A215").Value = "Montserrat"

# --- Row 4: Estados Unidos ---------------------------------------------
$ws.Range("B4").Value = 6788640
$ws.Range("C4").Value = 493
$ws.Range("E4").Value = 2519719
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 200217

# --- Row 25: Alemania ----------------------------------------------------
$ws.Range("B25").Value = 265142
$ws.Range("C25").Value = 298
$ws.Range("E25").Value = 16597

# --- Row 38: Kuwait ------------------------------------------------------
$ws.Range("B38").Value = 96999
$ws.Range("C38").Value = 698
$ws.Range("D38").Value = 87187
$ws.Range("E38").Value = 9241
$ws.Range("G38").Value = 3
$ws.Range("H38").Value = 571

# --- Row 42: Suecia -------------------------------------------------------
$ws.Range("B42").Value = 87575
$ws.Range("G42").Value = 5
$ws.Range("H42").Value = 5860

# --- Row 43: Paises Bajos (after swap) -----------------------------------
$ws.Range("B43").Value = 86320
$ws.Range("C43").Value = 1542
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 6260

# --- Row 44: China (after swap) -------------------------------------------
$ws.Range("B44").Value = 85214
$ws.Range("C44").Value = 12
$ws.Range("D44").Value = 80437
$ws.Range("E44").Value = 143
$ws.Range("H44").Value = 4634

# --- Row 61: Suiza ---------------------------------------------------------
$ws.Range("E61").Value = 6336
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 2029

# --- Row 68: Azerbaiyan ---------------------------------------------------
$ws.Range("B68").Value = 38658
$ws.Range("C68").Value = 141
$ws.Range("D68").Value = 36149
$ws.Range("E68").Value = 1940
$ws.Range("G68").Value = 3
$ws.Range("H68").Value = 569

# --- Row 114: Uganda (after swap) -----------------------------------------
$ws.Range("B114").Value = 5266
$ws.Range("C114").Value = 143
$ws.Range("D114").Value = 2404
$ws.Range("E114").Value = 2802
$ws.Range("G114").Value = 2
$ws.Range("H114").Value = 60

# --- Row 115: Suazilandia (after swap) ------------------------------------
$ws.Range("B115").Value = 5128
$ws.Range("D115").Value = 4401
$ws.Range("E115").Value = 626
$ws.Range("H115").Value = 101

# --- Row 150: Islandia ------------------------------------------------------
$ws.Range("B150").Value = 2189
$ws.Range("C150").Value = 15
$ws.Range("D150").Value = 2104
$ws.Range("E150").Value = 75

# --- Row 214: Islas Malvinas (after swap) ---------------------------------
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# --- Row 215: Montserrat (after swap) -------------------------------------
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
